# The sheet had an extra leading column A (row numbers 3/4/9/12 with a
# border style) that is no longer wanted. Deleting the entire column shifts
# B:F left into A:E, matching the target layout (dimension A1:E5).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").EntireColumn.Delete()

# The header that used to read "MODEL_CONDITION" (now shifted from E1 to D1)
# is retyped without the underscore.
$ws.Range("D1").Value = "MODELCONDITION"
